$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at X (column 24), shifting the existing "Comments" column
# from X to Y.
$ws.Columns.Item(24).Insert()

# Header for the newly inserted column.
$ws.Range("X1").Value = "lab"

# Fill in the "lab" values for the data rows: rows 2-123 were collected by
# Hakai, rows 124-153 by DFO.
$ws.Range("X2:X123").Value = "Hakai"
$ws.Range("X124:X153").Value = "DFO"

# Update the view state to match where the author was working.
$ws.Range("X124:X153").Select()
$excel.ActiveWindow.ScrollRow = 127
